# Slide 6 ("4-ii. Perform differential expression analysis"), the
# "Content Placeholder 6" body text box: split several runs apart to
# remove references to the old "tumor/normal" sample names (now
# "UHR"/"HBR") and update the library count, per the commit
# "removing more references to old sample names".

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(6)
$tf = $s.Shapes.Item(2).TextFrame
$tr = $tf.TextRange

# --- Paragraph 1: "In this step we will use cuffmerge and cuffdiff to:"
$para1 = $tr.Paragraphs(1)
$r = $para1.Runs(1)
$r.Text = "In this step we will use "
$r = $r.InsertAfter("cuffmerge")
$r = $r.InsertAfter(" and ")
$r = $r.InsertAfter("cuffdiff")
$r = $r.InsertAfter(" to:")

# --- Paragraph 2: "Combine expression estimates from our 4 libraries into more convenient files"
$para2 = $tr.Paragraphs(2)
$r = $para2.Runs(1)
$r.Text = "Combine expression estimates from our "
$r = $r.InsertAfter("6 ")
$r = $r.InsertAfter("libraries into more convenient files")

# --- Paragraph 3: "Combine expression estimates across replicates" (text unchanged)

# --- Paragraph 4: "Compare tumor vs. normal and identify significantly differentially expressed genes and isoforms (transcripts)"
$para4 = $tr.Paragraphs(4)
$r = $para4.Runs(1)
$r.Text = "Compare "
$r = $r.InsertAfter("UHR vs")
$r = $r.InsertAfter(". ")
$r = $r.InsertAfter("HBR ")
$r = $r.InsertAfter("and identify significantly differentially expressed genes and isoforms (transcripts)")

# --- Paragraph 5: "Note that these commands can get quite complicated when you have replicates" (text unchanged)

# --- Paragraph 6: "The positioning of spaces and commas, and grouping of libraries matters!" (text unchanged)

# --- Paragraph 7: "Comparisons" (text unchanged)

# --- Paragraph 8: "Compare Tumor vs. Normal using all replicates, for known (reference only mode) transcripts"
$para8 = $tr.Paragraphs(8)
$r = $para8.Runs(1)
$r.Text = "Compare "
$r = $r.InsertAfter("UHR vs")
$r = $r.InsertAfter(". ")
$r = $r.InsertAfter("HBR ")
$r = $r.InsertAfter("using all replicates, for known (reference only mode) transcripts")
